$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 731328207344.1487
$ws.Range("C3").Value = 243848336462.4445
$ws.Range("C4").Value = 38897308613.14926
$ws.Range("C5").Value = 33031459115.72026
$ws.Range("C6").Value = 24108716716.17286
$ws.Range("C7").Value = 13546617417.9198
$ws.Range("C8").Value = 11110165998.26605
$ws.Range("C9").Value = 8995956627.739824
$ws.Range("C10").Value = 8130294178.05073
$ws.Range("C11").Value = 8031994567.728826
$ws.Range("C12").Value = 7546205033.673535
$ws.Range("C13").Value = 7435826137.515828
$ws.Range("C14").Value = 6637440630.684655
$ws.Range("C15").Value = 6086459300.24696
$ws.Range("C16").Value = 5172735273.081258
$ws.Range("C17").Value = 4997271780.234427
$ws.Range("C18").Value = 4449279246.849122
$ws.Range("C19").Value = 3658360719.434418
$ws.Range("C20").Value = 3453039002.858096
$ws.Range("C21").Value = 3393770849.066844
$ws.Range("C22").Value = 3358900876.520487
$ws.Range("C23").Value = 3005546425.306632
$ws.Range("C24").Value = 2856450754.591782
$ws.Range("C25").Value = 2820653867.299805
$ws.Range("C26").Value = 2421578972.477175
